# Auto-generated Excel COM-interop script applying scheduled market-price
# refresh values to the Lich_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (48 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1404.4445
$ws.Range("I92").Value = 866.8570999999999
$ws.Range("J92").Value = 3286
$ws.Range("K92").Value = 866.8570999999999
$ws.Range("L92").Value = 3286
$ws.Range("M92").Value = 381.1429000000001
$ws.Range("N92").Value = -5782
$ws.Range("H97").Value = 3500
$ws.Range("J97").Value = 3500
$ws.Range("L97").Value = 10500
$ws.Range("N97").Value = -11492
$ws.Range("H100").Value = 1694.5294
$ws.Range("I100").Value = 1139.4166
$ws.Range("J100").Value = 3026.8
$ws.Range("K100").Value = 1139.4166
$ws.Range("L100").Value = 3026.8
$ws.Range("M100").Value = -598.4166
$ws.Range("N100").Value = -4108.8
$ws.Range("H101").Value = 2841804
$ws.Range("J101").Value = 1044.2858
$ws.Range("L101").Value = 3132.8574
$ws.Range("N101").Value = -6376.857400000001
$ws.Range("H113").Value = 4091.25
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H125").Value = 7258.0835
$ws.Range("I125").Value = 2435.9
$ws.Range("K125").Value = 21923.1
$ws.Range("M125").Value = -19463.1
$ws.Range("H131").Value = 52727.215
$ws.Range("I131").Value = 2756.3333
$ws.Range("J131").Value = 352552.5
$ws.Range("K131").Value = 8268.999899999999
$ws.Range("L131").Value = 1057657.5
$ws.Range("M131").Value = -3228.999899999999
$ws.Range("N131").Value = -1067737.5
$ws.Range("H132").Value = 1879.1489
$ws.Range("I132").Value = 1735.8292
$ws.Range("J132").Value = 2858.5
$ws.Range("K132").Value = 5207.487599999999
$ws.Range("L132").Value = 8575.5
$ws.Range("M132").Value = -2677.487599999999
$ws.Range("N132").Value = -13635.5
$ws.Range("H137").Value = 8987.799999999999
$ws.Range("I137").Value = 11938.3
$ws.Range("K137").Value = 35814.89999999999
$ws.Range("M137").Value = -33264.89999999999

# --- Sheet: ARM (16 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17979.309
$ws.Range("I32").Value = 16623.527
$ws.Range("K32").Value = 16623.527
$ws.Range("M32").Value = -16336.527
$ws.Range("H45").Value = 2488.889
$ws.Range("J45").Value = 2512.5
$ws.Range("L45").Value = 2512.5
$ws.Range("N45").Value = -3266.5
$ws.Range("H110").Value = 10025.833
$ws.Range("I110").Value = 10012.823
$ws.Range("K110").Value = 10012.823
$ws.Range("M110").Value = -7967.823
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- Sheet: BSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1234
$ws.Range("J80").Value = 1014
$ws.Range("L80").Value = 1014
$ws.Range("N80").Value = -3010
$ws.Range("H83").Value = 1234
$ws.Range("J83").Value = 1014
$ws.Range("L83").Value = 5070
$ws.Range("N83").Value = -15054

# --- Sheet: CRP (32 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 109999.5
$ws.Range("J20").Value = 109999.5
$ws.Range("L20").Value = 109999.5
$ws.Range("N20").Value = -110471.5
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H30").Value = 109999.5
$ws.Range("J30").Value = 109999.5
$ws.Range("L30").Value = 109999.5
$ws.Range("N30").Value = -110181.5
$ws.Range("H31").Value = 19952.387
$ws.Range("I31").Value = 1483
$ws.Range("K31").Value = 1483
$ws.Range("M31").Value = -1188
$ws.Range("H34").Value = 19952.387
$ws.Range("I34").Value = 1483
$ws.Range("K34").Value = 1483
$ws.Range("M34").Value = -1281
$ws.Range("H86").Value = 4713.7144
$ws.Range("J86").Value = 12139
$ws.Range("L86").Value = 12139
$ws.Range("N86").Value = -14385
$ws.Range("H89").Value = 4713.7144
$ws.Range("J89").Value = 12139
$ws.Range("L89").Value = 60695
$ws.Range("N89").Value = -71927
$ws.Range("H128").Value = 109999.5
$ws.Range("J128").Value = 109999.5
$ws.Range("L128").Value = 109999.5
$ws.Range("N128").Value = -119959.5

# --- Sheet: CUL (15 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 143580
$ws.Range("I7").Value = 500042.5
$ws.Range("J7").Value = 995
$ws.Range("K7").Value = 1500127.5
$ws.Range("L7").Value = 2985
$ws.Range("M7").Value = -1500015.5
$ws.Range("N7").Value = -3209
$ws.Range("H8").Value = 206.14285
$ws.Range("I8").Value = 206.14285
$ws.Range("K8").Value = 618.4285500000001
$ws.Range("M8").Value = -479.4285500000001
$ws.Range("H129").Value = 44975270
$ws.Range("J129").Value = 11113918
$ws.Range("L129").Value = 33341754
$ws.Range("N129").Value = -33351754

# --- Sheet: GSM (4 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1163.3334
$ws.Range("I3").Value = 1830
$ws.Range("K3").Value = 1830
$ws.Range("M3").Value = -1714

# --- Sheet: LTW (16 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2874.8
$ws.Range("I16").Value = 1704.0714
$ws.Range("K16").Value = 1704.0714
$ws.Range("M16").Value = -1534.0714
$ws.Range("H22").Value = 1871.091
$ws.Range("I22").Value = 847.4
$ws.Range("K22").Value = 847.4
$ws.Range("M22").Value = -552.4
$ws.Range("H27").Value = 1871.091
$ws.Range("I27").Value = 847.4
$ws.Range("K27").Value = 847.4
$ws.Range("M27").Value = -740.4
$ws.Range("H93").Value = 2934.0833
$ws.Range("I93").Value = 2933.2
$ws.Range("K93").Value = 2933.2
$ws.Range("M93").Value = -1685.2

# --- Sheet: WVR (24 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11629.667
$ws.Range("J62").Value = 11596
$ws.Range("L62").Value = 11596
$ws.Range("N62").Value = -12844
$ws.Range("H65").Value = 11629.667
$ws.Range("J65").Value = 11596
$ws.Range("L65").Value = 57980
$ws.Range("N65").Value = -64220
$ws.Range("H107").Value = 553.3570999999999
$ws.Range("I107").Value = 550.6923
$ws.Range("K107").Value = 1652.0769
$ws.Range("M107").Value = 267.9231
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178
